# Update the "想去人数" (F column) values on the "展览", "演出" and "全部类型"
# sheets to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1406
$ws1.Range("F5").Value  = 841
$ws1.Range("F6").Value  = 477
$ws1.Range("F8").Value  = 192
$ws1.Range("F10").Value = 50
$ws1.Range("F12").Value = 116
$ws1.Range("F13").Value = 1618
$ws1.Range("F14").Value = 205
$ws1.Range("F15").Value = 34
$ws1.Range("F17").Value = 73
$ws1.Range("F19").Value = 110
$ws1.Range("F24").Value = 49
$ws1.Range("F25").Value = 1458
$ws1.Range("F26").Value = 176

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 200

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1406
$ws4.Range("F6").Value  = 841
$ws4.Range("F9").Value  = 477
$ws4.Range("F12").Value = 192
$ws4.Range("F14").Value = 50
$ws4.Range("F16").Value = 116
$ws4.Range("F17").Value = 1618
$ws4.Range("F18").Value = 200
$ws4.Range("F19").Value = 205
$ws4.Range("F20").Value = 34
$ws4.Range("F22").Value = 73
$ws4.Range("F25").Value = 110
$ws4.Range("F36").Value = 49
$ws4.Range("F37").Value = 1458
$ws4.Range("F38").Value = 176
